# Update room names to room numbers in the timetable.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Timetable")

$ws.Range("D7").Value  = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("F7").Value  = "Flute MasterClass`n(Room G19)"

$ws.Range("C11").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("D11").Value = "Rehearsal with pianist`n(Room G22)"

$ws.Range("B19").Value = "Private Lesson with Stephane RETY `n(Room G19)"
$ws.Range("E19").Value = "Private Lesson with Stephane RETY & pianist `n(Room G19)"
$ws.Range("F19").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B27").Value = "Ensemble `n(Room G15)"
$ws.Range("C27").Value = "Ensemble `n(Room G15)"
$ws.Range("D27").Value = "Ensemble `n(Room G15)"
$ws.Range("E27").Value = "Ensemble `n(Room G15)"
$ws.Range("F27").Value = "Ensemble `n(Room G15)"
